# Complete tests for Accounts API
# Fill in the "Status" column (D) with "Done" for the newly completed
# Accounts API test rows (19-24), matching the styling already used for
# the other completed rows (green fill).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$doneFillColor = 5296274  # RGB(146, 208, 80) -> matches existing "Done" style

19..24 | ForEach-Object {
    $cell = $ws.Range("D$_")
    $cell.Value = "Done"
    $cell.Interior.Color = $doneFillColor
}

# Update the view state to match where the user ended up working
$ws.Activate()
$ws.Range("E24").Select()
